$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Adding User and Auth APIs:
# Row 5 was "GET users/:id" (fetch user by id, with query param "id : ......").
# It becomes "GET users/getdetails" and no longer needs a query param.
$ws.Range("C5").Value2 = "users/getdetails"
$ws.Range("D5").Value2 = ""

# Row 6 was "PUT users/:id" -> becomes "PUT users/update/:id"
$ws.Range("C6").Value2 = "users/update/:id"

# Row 7 was "DELETE users/:id" -> becomes "DELETE users/delete/:id"
$ws.Range("C7").Value2 = "users/delete/:id"

# Update the active selection/view to match the author's final cursor position.
$ws.Range("C7").Select()
